$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab
$ws.Name = "Country Stats"

# New column headers
$ws.Range("E2").Value = "Mean Profits Made (US$)"

# Column widths: C & D share the same width as before (col C), E and F get new widths
$ws.Columns.Item(3).ColumnWidth = 17.7109375
$ws.Columns.Item(4).ColumnWidth = 17.7109375
$ws.Columns.Item(5).ColumnWidth = 25.7109375
$ws.Columns.Item(6).ColumnWidth = 18.140625

$meanConst = 0.810965085

for ($r = 3; $r -le 45; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=D$r/C$r"
    $ws.Cells.Item($r, 6).Value = $meanConst
    $ws.Cells.Item($r, 7).Formula = "=(E$r-F$r)^2"
}

# Totals row
$ws.Cells.Item(46, 5).Formula = "=D46/C46"
$ws.Cells.Item(46, 6).Value = $meanConst

# Standard deviation row
$ws.Cells.Item(48, 3).Value = "Standard Deviation"
$ws.Cells.Item(48, 4).Formula = "=SQRT(SUM(G3:G45)/43)"

# View adjustments
$ws.Range("G12").Select()

$ws.PageSetup.Orientation = 1
